$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-05-06 Monday" "2024-05-07 Tuesday"

Replace-Text "20×46=" "80×44="
Replace-Text "73×79=" "35×39="
Replace-Text "86×46=" "33×56="
Replace-Text "38×53=" "15×24="
Replace-Text "28×83=" "85×61="

Replace-Text "25×54=" "33×95="
Replace-Text "14×73=" "31×18="
Replace-Text "89×38=" "69×42="
Replace-Text "16×21=" "92×22="
Replace-Text "44×49=" "66×63="

Replace-Text "63×76=" "84×74="
Replace-Text "33×23=" "99×59="
Replace-Text "63×45=" "59×33="
Replace-Text "70×97=" "73×94="
Replace-Text "92×84=" "79×63="

Replace-Text "85×75=" "66×78="
Replace-Text "61×41=" "44×12="
Replace-Text "63×55=" "87×40="
Replace-Text "70×51=" "44×69="
Replace-Text "71×73=" "96×96="

Replace-Text "22×24=" "25×25="
Replace-Text "33×87=" "82×56="
Replace-Text "34×86=" "68×61="
Replace-Text "30×98=" "44×80="
Replace-Text "13×64=" "42×13="
